$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-29 (Generation 0-27): Fitness -> 7894
$ws.Range("C2:C29").Value = 7894

# Rows 30-252 (Generation 28-250): Fitness -> 7569
$ws.Range("C30:C252").Value = 7569
